$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.255.06"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.343.30"
$ws.Range("E3").Value = "  +7.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.36"
$ws.Range("E5").Value = "  +3.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.51"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.13"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.388"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.347.81"
$ws.Range("E10").Value = "  +7.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.798"
$ws.Range("E11").Value = "  -4.67%  "
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.211.22"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.73"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000248"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.956.85"
$ws.Range("E16").Value = "  +7.34%  "
$ws.Range("B17").Value = "Toncoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.53"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.337.86"
$ws.Range("E18").Value = "  +7.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.66"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "493.90"
$ws.Range("E21").Value = "  +10.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.94"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000211"
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.30"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.70"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.76"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.18"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.504.15"
$ws.Range("E28").Value = "  +7.16%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.182"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.243"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.124"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.37"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.74"
$ws.Range("E35").Value = "  +6.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.157"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.55"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "510.22"
$ws.Range("E38").Value = "  +6.80%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.77"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.452"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.28"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.29"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.51"
$ws.Range("E44").Value = "  -7.99%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  +13.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.81"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.96"
$ws.Range("E48").Value = "  +4.75%  "
$ws.Range("E49").Value = "  +4.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0335"
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.07"
$ws.Range("E51").Value = "  +2.93%  "
